$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-5 from 45204 to 45207
$ws.Range("C2").Value = 45207
$ws.Range("C3").Value = 45207
$ws.Range("C4").Value = 45207
$ws.Range("C5").Value = 45207
